# "Generate Report for handoff"
#
# The localization-status report is regenerated: the row describing
# 86b72352-c94c-4240-a72c-f6ce183b670c.md moves from the top of the
# table to the bottom (it is now ready for a new handoff round), and
# its status flips from "Handed back: in sync with en-US" to
# "Ready for handoff" with refreshed handoff timestamps. The other two
# rows (ffff2d99924b... and ffffff7d49d452...) shift up to take its
# place, keeping their own data untouched.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "ffff2d99924b-233c-4d18-be84-ddfba0888239.md"
$ov.Range("B2").Value = "Handed back: in sync with en-US"
$ov.Range("C2").Value = "Handed back: in sync with en-US"

$ov.Range("A3").Value = "ffffff7d49d452-69c8-4708-926b-b6e395f4cf89.md"
$ov.Range("B3").Value = "Handed back: in sync with en-US"
$ov.Range("C3").Value = "Handed back: in sync with en-US"

$ov.Range("A4").Value = "86b72352-c94c-4240-a72c-f6ce183b670c.md"
$ov.Range("B4").Value = "Ready for handoff"
$ov.Range("C4").Value = "Ready for handoff"

# Rebuild the hyperlinks on column A in the new row order (targets are
# unchanged, only which row shows which file moved).
$ov.Cells.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/4cd537726be2f7049a93a35a65fdc5721c14f0f5/e2e/ffff2d99924b-233c-4d18-be84-ddfba0888239.md", "", "", "ffff2d99924b-233c-4d18-be84-ddfba0888239.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/4cd537726be2f7049a93a35a65fdc5721c14f0f5/e2e/ffffff7d49d452-69c8-4708-926b-b6e395f4cf89.md", "", "", "ffffff7d49d452-69c8-4708-926b-b6e395f4cf89.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/4cd537726be2f7049a93a35a65fdc5721c14f0f5/e2e/86b72352-c94c-4240-a72c-f6ce183b670c.md", "", "", "86b72352-c94c-4240-a72c-f6ce183b670c.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/4cd537726be2f7049a93a35a65fdc5721c14f0f5/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "ffff2d99924b-233c-4d18-be84-ddfba0888239.md"
$zh.Range("B2").Value = "Handed back: in sync with en-US"
$zh.Range("C2").Value = "fe219963-43f6-4357-a768-730c515b9e56.465c21a7bc3af6829ddc588a5812abf7541dc467.zh-cn.xlf"
$zh.Range("D2").Value = "2016-01-20 03:51:39"
$zh.Range("E2").Value = "fe219963-43f6-4357-a768-730c515b9e56.md"
$zh.Range("F2").Value = "fe219963-43f6-4357-a768-730c515b9e56.465c21a7bc3af6829ddc588a5812abf7541dc467.zh-cn.xlf"
$zh.Range("G2").Value = "2016-01-20 03:52:22"
$zh.Range("H2").Value = "Include"

$zh.Range("A3").Value = "ffffff7d49d452-69c8-4708-926b-b6e395f4cf89.md"
$zh.Range("B3").Value = "Handed back: in sync with en-US"
$zh.Range("C3").Value = "fe219963-43f6-4357-a768-730c515b9e56.465c21a7bc3af6829ddc588a5812abf7541dc467.zh-cn.xlf"
$zh.Range("D3").Value = "2016-01-20 03:51:39"
$zh.Range("E3").Value = "fe219963-43f6-4357-a768-730c515b9e56.md"
$zh.Range("F3").Value = "fe219963-43f6-4357-a768-730c515b9e56.465c21a7bc3af6829ddc588a5812abf7541dc467.zh-cn.xlf"
$zh.Range("G3").Value = "2016-01-20 03:52:22"
$zh.Range("H3").Value = "Include"

$zh.Range("A4").Value = "86b72352-c94c-4240-a72c-f6ce183b670c.md"
$zh.Range("B4").Value = "Ready for handoff"
$zh.Range("C4").Value = "86b72352-c94c-4240-a72c-f6ce183b670c.bda3accd15fba864b15cb74fc19a4f35fd9362a9.zh-cn.xlf"
$zh.Range("D4").Value = "2016-01-20 03:55:35"
$zh.Range("E4").Value = "86b72352-c94c-4240-a72c-f6ce183b670c.md"
$zh.Range("F4").Value = "86b72352-c94c-4240-a72c-f6ce183b670c.bda3accd15fba864b15cb74fc19a4f35fd9362a9.zh-cn.xlf"
$zh.Range("G4").Value = "2016-01-20 03:54:33"
$zh.Range("H4").Value = "Include"

$zh.Cells.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/4cd537726be2f7049a93a35a65fdc5721c14f0f5/e2e/ffff2d99924b-233c-4d18-be84-ddfba0888239.md", "", "", "ffff2d99924b-233c-4d18-be84-ddfba0888239.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4eda289f1383f862ffae1bf4d8b754e66b85fc4e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/fe219963-43f6-4357-a768-730c515b9e56.465c21a7bc3af6829ddc588a5812abf7541dc467.zh-cn.xlf", "", "", "fe219963-43f6-4357-a768-730c515b9e56.465c21a7bc3af6829ddc588a5812abf7541dc467.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/51cca93bbbfe32f53252909e65bffc23773e4472/e2e/fe219963-43f6-4357-a768-730c515b9e56.md", "", "", "fe219963-43f6-4357-a768-730c515b9e56.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/04f4c3f484cf2362d2e5ec07b217c7a1e5730297/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/fe219963-43f6-4357-a768-730c515b9e56.465c21a7bc3af6829ddc588a5812abf7541dc467.zh-cn.xlf", "", "", "fe219963-43f6-4357-a768-730c515b9e56.465c21a7bc3af6829ddc588a5812abf7541dc467.zh-cn.xlf") | Out-Null

$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/4cd537726be2f7049a93a35a65fdc5721c14f0f5/e2e/ffffff7d49d452-69c8-4708-926b-b6e395f4cf89.md", "", "", "ffffff7d49d452-69c8-4708-926b-b6e395f4cf89.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4eda289f1383f862ffae1bf4d8b754e66b85fc4e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/fe219963-43f6-4357-a768-730c515b9e56.465c21a7bc3af6829ddc588a5812abf7541dc467.zh-cn.xlf", "", "", "fe219963-43f6-4357-a768-730c515b9e56.465c21a7bc3af6829ddc588a5812abf7541dc467.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/51cca93bbbfe32f53252909e65bffc23773e4472/e2e/fe219963-43f6-4357-a768-730c515b9e56.md", "", "", "fe219963-43f6-4357-a768-730c515b9e56.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/04f4c3f484cf2362d2e5ec07b217c7a1e5730297/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/fe219963-43f6-4357-a768-730c515b9e56.465c21a7bc3af6829ddc588a5812abf7541dc467.zh-cn.xlf", "", "", "fe219963-43f6-4357-a768-730c515b9e56.465c21a7bc3af6829ddc588a5812abf7541dc467.zh-cn.xlf") | Out-Null

$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/4cd537726be2f7049a93a35a65fdc5721c14f0f5/e2e/86b72352-c94c-4240-a72c-f6ce183b670c.md", "", "", "86b72352-c94c-4240-a72c-f6ce183b670c.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b85c76535210dc4338b86632ea6d1354cffd4cc9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/86b72352-c94c-4240-a72c-f6ce183b670c.bda3accd15fba864b15cb74fc19a4f35fd9362a9.zh-cn.xlf", "", "", "86b72352-c94c-4240-a72c-f6ce183b670c.bda3accd15fba864b15cb74fc19a4f35fd9362a9.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/42e7890d452f40a1ff513ba727ae4197da55a059/e2e/86b72352-c94c-4240-a72c-f6ce183b670c.md", "", "", "86b72352-c94c-4240-a72c-f6ce183b670c.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/005e7e5d227406a2c941330d99cd8d61996cd68b/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/86b72352-c94c-4240-a72c-f6ce183b670c.bda3accd15fba864b15cb74fc19a4f35fd9362a9.zh-cn.xlf", "", "", "86b72352-c94c-4240-a72c-f6ce183b670c.bda3accd15fba864b15cb74fc19a4f35fd9362a9.zh-cn.xlf") | Out-Null

$zh.Hyperlinks.Add($zh.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/4cd537726be2f7049a93a35a65fdc5721c14f0f5/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "ffff2d99924b-233c-4d18-be84-ddfba0888239.md"
$de.Range("B2").Value = "Handed back: in sync with en-US"
$de.Range("C2").Value = "fe219963-43f6-4357-a768-730c515b9e56.465c21a7bc3af6829ddc588a5812abf7541dc467.de-de.xlf"
$de.Range("D2").Value = "2016-01-20 03:51:49"
$de.Range("E2").Value = "fe219963-43f6-4357-a768-730c515b9e56.md"
$de.Range("F2").Value = "fe219963-43f6-4357-a768-730c515b9e56.465c21a7bc3af6829ddc588a5812abf7541dc467.de-de.xlf"
$de.Range("G2").Value = "2016-01-20 03:52:40"
$de.Range("H2").Value = "Include"

$de.Range("A3").Value = "ffffff7d49d452-69c8-4708-926b-b6e395f4cf89.md"
$de.Range("B3").Value = "Handed back: in sync with en-US"
$de.Range("C3").Value = "fe219963-43f6-4357-a768-730c515b9e56.465c21a7bc3af6829ddc588a5812abf7541dc467.de-de.xlf"
$de.Range("D3").Value = "2016-01-20 03:51:49"
$de.Range("E3").Value = "fe219963-43f6-4357-a768-730c515b9e56.md"
$de.Range("F3").Value = "fe219963-43f6-4357-a768-730c515b9e56.465c21a7bc3af6829ddc588a5812abf7541dc467.de-de.xlf"
$de.Range("G3").Value = "2016-01-20 03:52:40"
$de.Range("H3").Value = "Include"

$de.Range("A4").Value = "86b72352-c94c-4240-a72c-f6ce183b670c.md"
$de.Range("B4").Value = "Ready for handoff"
$de.Range("C4").Value = "86b72352-c94c-4240-a72c-f6ce183b670c.bda3accd15fba864b15cb74fc19a4f35fd9362a9.de-de.xlf"
$de.Range("D4").Value = "2016-01-20 03:55:45"
$de.Range("E4").Value = "86b72352-c94c-4240-a72c-f6ce183b670c.md"
$de.Range("F4").Value = "86b72352-c94c-4240-a72c-f6ce183b670c.bda3accd15fba864b15cb74fc19a4f35fd9362a9.de-de.xlf"
$de.Range("G4").Value = "2016-01-20 03:54:51"
$de.Range("H4").Value = "Include"

$de.Cells.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/4cd537726be2f7049a93a35a65fdc5721c14f0f5/e2e/ffff2d99924b-233c-4d18-be84-ddfba0888239.md", "", "", "ffff2d99924b-233c-4d18-be84-ddfba0888239.md") | Out-Null
$de.Hyperlinks.Add($de.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d1a335af176d7d3ebb9048e30619daf300b4188a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/fe219963-43f6-4357-a768-730c515b9e56.465c21a7bc3af6829ddc588a5812abf7541dc467.de-de.xlf", "", "", "fe219963-43f6-4357-a768-730c515b9e56.465c21a7bc3af6829ddc588a5812abf7541dc467.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/c6524b33452ae700e442bde2235f36db64ac0427/e2e/fe219963-43f6-4357-a768-730c515b9e56.md", "", "", "fe219963-43f6-4357-a768-730c515b9e56.md") | Out-Null
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/73b4dbd531619f988263456adbb9f61b11bb39c4/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/fe219963-43f6-4357-a768-730c515b9e56.465c21a7bc3af6829ddc588a5812abf7541dc467.de-de.xlf", "", "", "fe219963-43f6-4357-a768-730c515b9e56.465c21a7bc3af6829ddc588a5812abf7541dc467.de-de.xlf") | Out-Null

$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/4cd537726be2f7049a93a35a65fdc5721c14f0f5/e2e/ffffff7d49d452-69c8-4708-926b-b6e395f4cf89.md", "", "", "ffffff7d49d452-69c8-4708-926b-b6e395f4cf89.md") | Out-Null
$de.Hyperlinks.Add($de.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d1a335af176d7d3ebb9048e30619daf300b4188a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/fe219963-43f6-4357-a768-730c515b9e56.465c21a7bc3af6829ddc588a5812abf7541dc467.de-de.xlf", "", "", "fe219963-43f6-4357-a768-730c515b9e56.465c21a7bc3af6829ddc588a5812abf7541dc467.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/c6524b33452ae700e442bde2235f36db64ac0427/e2e/fe219963-43f6-4357-a768-730c515b9e56.md", "", "", "fe219963-43f6-4357-a768-730c515b9e56.md") | Out-Null
$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/73b4dbd531619f988263456adbb9f61b11bb39c4/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/fe219963-43f6-4357-a768-730c515b9e56.465c21a7bc3af6829ddc588a5812abf7541dc467.de-de.xlf", "", "", "fe219963-43f6-4357-a768-730c515b9e56.465c21a7bc3af6829ddc588a5812abf7541dc467.de-de.xlf") | Out-Null

$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/4cd537726be2f7049a93a35a65fdc5721c14f0f5/e2e/86b72352-c94c-4240-a72c-f6ce183b670c.md", "", "", "86b72352-c94c-4240-a72c-f6ce183b670c.md") | Out-Null
$de.Hyperlinks.Add($de.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/be6ecd672a39e82b605701d7145062bd0cb34bba/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/86b72352-c94c-4240-a72c-f6ce183b670c.bda3accd15fba864b15cb74fc19a4f35fd9362a9.de-de.xlf", "", "", "86b72352-c94c-4240-a72c-f6ce183b670c.bda3accd15fba864b15cb74fc19a4f35fd9362a9.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/32531a78009b88bed3f392f856de9b27bd34ce16/e2e/86b72352-c94c-4240-a72c-f6ce183b670c.md", "", "", "86b72352-c94c-4240-a72c-f6ce183b670c.md") | Out-Null
$de.Hyperlinks.Add($de.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/7ded0cf9da41885a22a842c4b11c41ab95c72109/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/86b72352-c94c-4240-a72c-f6ce183b670c.bda3accd15fba864b15cb74fc19a4f35fd9362a9.de-de.xlf", "", "", "86b72352-c94c-4240-a72c-f6ce183b670c.bda3accd15fba864b15cb74fc19a4f35fd9362a9.de-de.xlf") | Out-Null

$de.Hyperlinks.Add($de.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/4cd537726be2f7049a93a35a65fdc5721c14f0f5/.localization-config", "", "", ".localization-config") | Out-Null

Write-Output "done"
